# excel/new_ph_pl/new_ph_pl.xlsx :: "day" sheet
# 1) D117:D143 (bsecode) were stored as text; convert them to real numbers.
# 2) Append new rows 144-152 scraped from the "stock.yaml" break-out.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("day")

# --- 1. Fix bsecode column (D) for existing rows 117-143: text -> number ---
$bsecodes = @{
    117 = 500530
    118 = 532977
    119 = 508869
    120 = 532644
    121 = 539448
    122 = 523642
    123 = 500420
    124 = 542650
    125 = 532281
    126 = 500087
    127 = 524804
    128 = 532215
    129 = 542830
    130 = 500570
    131 = 532868
    132 = 500260
    133 = 500670
    134 = 533278
    135 = 540222
    136 = 500400
    137 = 513599
    138 = 500312
    139 = 532234
    140 = 530965
    141 = 500113
    142 = 532461
    143 = 532754
}

foreach ($row in $bsecodes.Keys) {
    $ws.Cells.Item($row, 4).Value = $bsecodes[$row]
}

# --- 2. Append the new data rows (144-152) ---
# columns: sr, nsecode, name, bsecode, per_chg, close, volume, timeframe, Date Time
$newRows = @(
    @(144, 1, "MCX",        "Multi Commodity Exchange Of India Limited", "534091", 0.32,  3771.35, 549537,    "day", "12/07/2024 11:41:53"),
    @(145, 2, "BAJAJFINSV", "Bajaj Finserv Limited",                     "532978", 0.59,  1597.25, 1003699,   "day", "12/07/2024 11:41:53"),
    @(146, 3, "BATAINDIA",  "Bata India Limited",                        "500043", -1.34, 1514,    332578,    "day", "12/07/2024 11:41:53"),
    @(147, 4, "BHARTIARTL", "Bharti Airtel Limited",                     "532454", -0.33, 1433.25, 6342136,   "day", "12/07/2024 11:41:53"),
    @(148, 5, "IPCALAB",    "Ipca Laboratories Limited",                 "524494", 0.09,  1225,    240196,    "day", "12/07/2024 11:41:53"),
    @(149, 6, "NTPC",       "Ntpc Limited",                              "532555", 0.01,  377.15,  16151583,  "day", "12/07/2024 11:41:53"),
    @(150, 7, "INDIACEM",   "The India Cements Limited",                 "530005", -2.31, 297.85,  4144724,   "day", "12/07/2024 11:41:53"),
    @(151, 8, "RBLBANK",    "Rbl Bank Limited",                          "540065", 0.26,  246.05,  5256242,   "day", "12/07/2024 11:41:53"),
    @(152, 9, "IDEA",       "Idea Cellular Limited",                     "532822", -2.84, 16.09,   630743126, "day", "12/07/2024 11:41:53")
)

foreach ($r in $newRows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    # bsecode kept as text (matches source: inline string, not numeric) -
    # leading apostrophe forces text entry without leaving a stray format.
    $ws.Cells.Item($rowNum, 4).Value = "'" + $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]
    $ws.Cells.Item($rowNum, 7).Value = $r[7]
    $ws.Cells.Item($rowNum, 8).Value = $r[8]
    $ws.Cells.Item($rowNum, 9).Value = $r[9]
}
